$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3857.5
$ws.Range("I74").Value = 3320.75
$ws.Range("J74").Value = 4036.4167
$ws.Range("K74").Value = 3320.75
$ws.Range("L74").Value = 4036.4167
$ws.Range("M74").Value = -2384.75
$ws.Range("N74").Value = -5908.4167
$ws.Range("H77").Value = 3857.5
$ws.Range("I77").Value = 3320.75
$ws.Range("J77").Value = 4036.4167
$ws.Range("K77").Value = 16603.75
$ws.Range("L77").Value = 20182.0835
$ws.Range("M77").Value = -11923.75
$ws.Range("N77").Value = -29542.0835
$ws.Range("H80").Value = 2406.8333
$ws.Range("I80").Value = 744.8095
$ws.Range("J80").Value = 6284.8887
$ws.Range("K80").Value = 2234.4285
$ws.Range("L80").Value = 18854.6661
$ws.Range("M80").Value = -1236.4285
$ws.Range("N80").Value = -20850.6661
$ws.Range("H82").Value = 3870.2856
$ws.Range("J82").Value = 9400
$ws.Range("L82").Value = 28200
$ws.Range("N82").Value = -29012
$ws.Range("H83").Value = 2406.8333
$ws.Range("I83").Value = 744.8095
$ws.Range("J83").Value = 6284.8887
$ws.Range("K83").Value = 6703.2855
$ws.Range("L83").Value = 56563.99830000001
$ws.Range("M83").Value = -1711.2855
$ws.Range("N83").Value = -66547.99830000001
$ws.Range("H85").Value = 3870.2856
$ws.Range("J85").Value = 9400
$ws.Range("L85").Value = 28200
$ws.Range("N85").Value = -31008
$ws.Range("H96").Value = 1607.28
$ws.Range("I96").Value = 506.4375
$ws.Range("J96").Value = 3564.3333
$ws.Range("K96").Value = 1519.3125
$ws.Range("L96").Value = 10692.9999
$ws.Range("M96").Value = -146.3125
$ws.Range("N96").Value = -13438.9999
$ws.Range("H113").Value = 9897
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 11867.571
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 11867.571
$ws.Range("M113").Value = 254
$ws.Range("N113").Value = -18375.571
$ws.Range("H129").Value = 952.9167
$ws.Range("I129").Value = 235.1
$ws.Range("J129").Value = 1096.48
$ws.Range("K129").Value = 705.3
$ws.Range("L129").Value = 3289.44
$ws.Range("M129").Value = 4294.7
$ws.Range("N129").Value = -13289.44
$ws.Range("H137").Value = 1332.3529
$ws.Range("I137").Value = 1025.2
$ws.Range("J137").Value = 1771.1428
$ws.Range("K137").Value = 3075.6
$ws.Range("L137").Value = 5313.428400000001
$ws.Range("M137").Value = -525.6000000000004
$ws.Range("N137").Value = -10413.4284
$ws.Range("H139").Value = 70178.336
$ws.Range("J139").Value = 70178.336
$ws.Range("L139").Value = 70178.336
$ws.Range("N139").Value = -80458.336

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 431734.38
$ws.Range("I32").Value = 5232.1875
$ws.Range("J32").Value = 9530448
$ws.Range("K32").Value = 5232.1875
$ws.Range("L32").Value = 9530448
$ws.Range("M32").Value = -4945.1875
$ws.Range("N32").Value = -9531022
$ws.Range("H37").Value = 142862900
$ws.Range("J37").Value = 8100.6665
$ws.Range("L37").Value = 8100.6665
$ws.Range("N37").Value = -8646.666499999999
$ws.Range("H45").Value = 3630.8845
$ws.Range("I45").Value = 3824.8572
$ws.Range("K45").Value = 3824.8572
$ws.Range("M45").Value = -3447.8572
$ws.Range("H63").Value = 5916.364
$ws.Range("I63").Value = 3180
$ws.Range("J63").Value = 9200
$ws.Range("K63").Value = 3180
$ws.Range("L63").Value = 9200
$ws.Range("M63").Value = -2494
$ws.Range("N63").Value = -10572
$ws.Range("H66").Value = 5916.364
$ws.Range("I66").Value = 3180
$ws.Range("J66").Value = 9200
$ws.Range("K66").Value = 15900
$ws.Range("L66").Value = 46000
$ws.Range("M66").Value = -12468
$ws.Range("N66").Value = -52864
$ws.Range("H74").Value = 1052.2667
$ws.Range("I74").Value = 746
$ws.Range("J74").Value = 1256.4445
$ws.Range("K74").Value = 746
$ws.Range("L74").Value = 1256.4445
$ws.Range("M74").Value = 128
$ws.Range("N74").Value = -3004.4445
$ws.Range("H77").Value = 1052.2667
$ws.Range("I77").Value = 746
$ws.Range("J77").Value = 1256.4445
$ws.Range("K77").Value = 3730
$ws.Range("L77").Value = 6282.2225
$ws.Range("M77").Value = 638
$ws.Range("N77").Value = -15018.2225
$ws.Range("H88").Value = 5486.6665
$ws.Range("I88").Value = 2133.3333
$ws.Range("J88").Value = 7163.3335
$ws.Range("K88").Value = 2133.3333
$ws.Range("L88").Value = 7163.3335
$ws.Range("M88").Value = -1727.3333
$ws.Range("N88").Value = -7975.3335
$ws.Range("H91").Value = 5486.6665
$ws.Range("I91").Value = 2133.3333
$ws.Range("J91").Value = 7163.3335
$ws.Range("K91").Value = 2133.3333
$ws.Range("L91").Value = 7163.3335
$ws.Range("M91").Value = -729.3332999999998
$ws.Range("N91").Value = -9971.333500000001
$ws.Range("H122").Value = 19324.924
$ws.Range("I122").Value = 24021.9
$ws.Range("J122").Value = 3668.3333
$ws.Range("K122").Value = 72065.70000000001
$ws.Range("L122").Value = 11004.9999
$ws.Range("M122").Value = -69615.70000000001
$ws.Range("N122").Value = -15904.9999
$ws.Range("H138").Value = 63528.57
$ws.Range("J138").Value = 63528.57
$ws.Range("L138").Value = 63528.57
$ws.Range("N138").Value = -73808.57000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1158.8889
$ws.Range("I99").Value = 1158.8889
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1158.8889
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 339.1111000000001
$ws.Range("N99").ClearContents()
$ws.Range("H138").Value = 69740
$ws.Range("J138").Value = 69740
$ws.Range("L138").Value = 69740
$ws.Range("N138").Value = -80020

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2688.12
$ws.Range("I31").Value = 1576.0889
$ws.Range("J31").Value = 3597.9636
$ws.Range("K31").Value = 1576.0889
$ws.Range("L31").Value = 3597.9636
$ws.Range("M31").Value = -1281.0889
$ws.Range("N31").Value = -4187.9636
$ws.Range("H34").Value = 2688.12
$ws.Range("I34").Value = 1576.0889
$ws.Range("J34").Value = 3597.9636
$ws.Range("K34").Value = 1576.0889
$ws.Range("L34").Value = 3597.9636
$ws.Range("M34").Value = -1374.0889
$ws.Range("N34").Value = -4001.9636
$ws.Range("H60").Value = 7924
$ws.Range("J60").Value = 8201
$ws.Range("L60").Value = 8201
$ws.Range("N60").Value = -9223
$ws.Range("H122").Value = 2327
$ws.Range("I122").Value = 1997
$ws.Range("J122").Value = 2987
$ws.Range("K122").Value = 5991
$ws.Range("L122").Value = 8961
$ws.Range("M122").Value = -3541
$ws.Range("N122").Value = -13861
$ws.Range("H138").Value = 48200
$ws.Range("J138").Value = 48200
$ws.Range("L138").Value = 48200
$ws.Range("N138").Value = -58480
$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360
$ws.Range("H141").Value = 23165.334
$ws.Range("J141").Value = 29600
$ws.Range("L141").Value = 29600
$ws.Range("N141").Value = -39960

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 469.45715
$ws.Range("I122").Value = 342.48276
$ws.Range("J122").Value = 1083.1666
$ws.Range("K122").Value = 3082.34484
$ws.Range("L122").Value = 9748.499400000001
$ws.Range("M122").Value = -632.3448399999997
$ws.Range("N122").Value = -14648.4994

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2477.889
$ws.Range("I122").Value = 2372.182
$ws.Range("J122").Value = 2644
$ws.Range("K122").Value = 7116.545999999999
$ws.Range("L122").Value = 7932
$ws.Range("M122").Value = -4666.545999999999
$ws.Range("N122").Value = -12832
$ws.Range("H138").Value = 68400
$ws.Range("J138").Value = 68400
$ws.Range("L138").Value = 68400
$ws.Range("N138").Value = -78680
$ws.Range("H140").Value = 89894.5
$ws.Range("J140").Value = 89894.5
$ws.Range("L140").Value = 89894.5
$ws.Range("N140").Value = -100254.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H122").Value = 3369.2856
$ws.Range("I122").Value = 2875.375
$ws.Range("J122").Value = 4027.8333
$ws.Range("K122").Value = 8626.125
$ws.Range("L122").Value = 12083.4999
$ws.Range("M122").Value = -6176.125
$ws.Range("N122").Value = -16983.4999
$ws.Range("H132").Value = 4190
$ws.Range("I132").Value = 6059.0835
$ws.Range("K132").Value = 18177.2505
$ws.Range("M132").Value = -15647.2505
$ws.Range("H139").Value = 79750
$ws.Range("J139").Value = 79750
$ws.Range("L139").Value = 79750
$ws.Range("N139").Value = -90030

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6500
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 6500
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 6500
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -7540
$ws.Range("H100").Value = 1275.8889
$ws.Range("I100").Value = 1185.375
$ws.Range("K100").Value = 2370.75
$ws.Range("M100").Value = -1829.75
$ws.Range("H122").Value = 2724
$ws.Range("I122").Value = 2632.8572
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 7898.571599999999
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -5448.571599999999
$ws.Range("N122").Value = -16900
